# Updated viz for Day 4.
# The "2019" worksheet tracks Advent of Code daily puzzle times.
# Day 4 had only a placeholder title ("Day 4: TITLE") with no data yet;
# now that the puzzle is known/solved, fill in its real title and times.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2019")
$ws.Select()

# Row 8 corresponds to "Day 4". Replace the placeholder title with the
# real puzzle name, and fill in the recorded times (start->part1, part1->part2
# differences are computed by existing formulas in columns D and G).
$ws.Range("B8").Value = "Day 4: Secure Container"
$ws.Range("C8").Value = 0.00337962962962963
$ws.Range("E8").Value = 0.02245370370370371
$ws.Range("F8").Value = 0.005046296296296296
$ws.Range("H8").Value = "11th"

# Leave the cursor where the author left it after editing.
$ws.Range("G37").Select()
